# FantaSPL_Classifica: switch the ranking table so column A shows the
# numeric Rank (used to link to each team), column B the Fantasy Team
# Name and column C the Total Points Scored. The old "Player Name" /
# "Rank" columns are dropped (the rank now lives in column A), so the
# sheet shrinks from A1:D24 to A1:C24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Rank"
$ws.Cells.Item(1, 2).Value = "Fantasy Team Name"
$ws.Cells.Item(1, 3).Value = "Total Points Scored"

# --- Data rows (Rank, Fantasy Team Name, Total Points Scored) -------
$data = @(
    @(1,  "Limonta United", 67),
    @(1,  "Football Meta Academy", 67),
    @(1,  "SPL:Solo Per Letette (ciccu)", 67),
    @(4,  "Caledonians", 65),
    @(5,  "Omanta", 63),
    @(6,  "T'eamCulo", 61),
    @(6,  "CHIAVO VERONA", 61),
    @(8,  "Non è la seconda squadra di Mazzu, è la prima", 59),
    @(9,  "Aldo Ritmo", 58),
    @(9,  "LA PALLA NON ERA USCITA", 58),
    @(9,  "Rapid Viennetta", 58),
    @(12, "Mazzu è ok", 54),
    @(13, "Cesarino’s", 52),
    @(14, "Mazzu doveva Vincere", 51),
    @(15, "BARBA FC", 50),
    @(16, "Multiple Cancers", 47),
    @(17, "Latin Supremacy", 45),
    @(18, "Beverly INPS", 42),
    @(19, "Affori Grizzlies", 39),
    @(19, "I nemici di mazzu", 39),
    @(21, "Fury Bonds", 35),
    @(22, "Fel Lazio", 33),
    @(23, "Artificially Degenerated", 9)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Drop the now-unused column D (old Rank column) ------------------
$ws.Columns.Item(4).Delete()
